$d = $word.ActiveDocument

$replacements = @(
    @{old="31×90=2790"; new="52×94=4888"},
    @{old="73×13=949"; new="76×58=4408"},
    @{old="71×26=1846"; new="62×82=5084"},
    @{old="83×80=6640"; new="25×42=1050"},
    @{old="79×14=1106"; new="97×47=4559"},
    @{old="19×71=1349"; new="49×57=2793"},
    @{old="88×97=8536"; new="34×36=1224"},
    @{old="12×86=1032"; new="20×30=600"},
    @{old="64×87=5568"; new="29×46=1334"},
    @{old="14×49=686"; new="38×14=532"},
    @{old="81×30=2430"; new="88×61=5368"},
    @{old="56×36=2016"; new="65×89=5785"},
    @{old="11×32=352"; new="74×68=5032"},
    @{old="90×76=6840"; new="81×98=7938"},
    @{old="68×47=3196"; new="64×96=6144"},
    @{old="70×33=2310"; new="80×71=5680"},
    @{old="26×82=2132"; new="45×79=3555"},
    @{old="45×83=3735"; new="23×93=2139"},
    @{old="45×21=945"; new="37×92=3404"},
    @{old="25×92=2300"; new="21×33=693"},
    @{old="60×30=1800"; new="41×97=3977"},
    @{old="32×71=2272"; new="56×79=4424"},
    @{old="81×59=4779"; new="77×65=5005"},
    @{old="34×16=544"; new="30×38=1140"},
    @{old="26×17=442"; new="43×88=3784"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
